$d = $word.ActiveDocument

$d.Content.Find.Execute("57-8=", $false, $false, $false, $false, $false, $true, 1, $false, "73-55=", 2) | Out-Null
$d.Content.Find.Execute("41+28=", $false, $false, $false, $false, $false, $true, 1, $false, "72+24=", 2) | Out-Null
$d.Content.Find.Execute("81-9=", $false, $false, $false, $false, $false, $true, 1, $false, "93-86=", 2) | Out-Null
$d.Content.Find.Execute("41+33=", $false, $false, $false, $false, $false, $true, 1, $false, "90-48=", 2) | Out-Null
$d.Content.Find.Execute("97-30=", $false, $false, $false, $false, $false, $true, 1, $false, "95-1=", 2) | Out-Null
$d.Content.Find.Execute("88-17=", $false, $false, $false, $false, $false, $true, 1, $false, "14-4=", 2) | Out-Null
$d.Content.Find.Execute("3+11=", $false, $false, $false, $false, $false, $true, 1, $false, "38+49=", 2) | Out-Null
$d.Content.Find.Execute("49-30=", $false, $false, $false, $false, $false, $true, 1, $false, "94-92=", 2) | Out-Null
$d.Content.Find.Execute("11+2=", $false, $false, $false, $false, $false, $true, 1, $false, "58+8=", 2) | Out-Null
$d.Content.Find.Execute("54-27=", $false, $false, $false, $false, $false, $true, 1, $false, "20-0=", 2) | Out-Null
$d.Content.Find.Execute("7+84=", $false, $false, $false, $false, $false, $true, 1, $false, "92-8=", 2) | Out-Null
$d.Content.Find.Execute("84+9=", $false, $false, $false, $false, $false, $true, 1, $false, "83-47=", 2) | Out-Null
$d.Content.Find.Execute("3+23=", $false, $false, $false, $false, $false, $true, 1, $false, "2+0=", 2) | Out-Null
$d.Content.Find.Execute("48-2=", $false, $false, $false, $false, $false, $true, 1, $false, "37-14=", 2) | Out-Null
$d.Content.Find.Execute("73-72=", $false, $false, $false, $false, $false, $true, 1, $false, "7+25=", 2) | Out-Null
$d.Content.Find.Execute("35+38=", $false, $false, $false, $false, $false, $true, 1, $false, "96-62=", 2) | Out-Null
$d.Content.Find.Execute("62+32=", $false, $false, $false, $false, $false, $true, 1, $false, "38+42=", 2) | Out-Null
$d.Content.Find.Execute("47+25=", $false, $false, $false, $false, $false, $true, 1, $false, "39+42=", 2) | Out-Null
$d.Content.Find.Execute("21+59=", $false, $false, $false, $false, $false, $true, 1, $false, "67-4=", 2) | Out-Null
$d.Content.Find.Execute("46+38=", $false, $false, $false, $false, $false, $true, 1, $false, "27-23=", 2) | Out-Null
$d.Content.Find.Execute("74+20=", $false, $false, $false, $false, $false, $true, 1, $false, "9+38=", 2) | Out-Null
$d.Content.Find.Execute("76+2=", $false, $false, $false, $false, $false, $true, 1, $false, "99-33=", 2) | Out-Null
$d.Content.Find.Execute("42-33=", $false, $false, $false, $false, $false, $true, 1, $false, "73-40=", 2) | Out-Null
$d.Content.Find.Execute("26+40=", $false, $false, $false, $false, $false, $true, 1, $false, "52-40=", 2) | Out-Null
$d.Content.Find.Execute("86-46=", $false, $false, $false, $false, $false, $true, 1, $false, "31+19=", 2) | Out-Null
$d.Content.Find.Execute("17+26=", $false, $false, $false, $false, $false, $true, 1, $false, "51+38=", 2) | Out-Null
$d.Content.Find.Execute("17+10=", $false, $false, $false, $false, $false, $true, 1, $false, "58-13=", 2) | Out-Null
$d.Content.Find.Execute("36+56=", $false, $false, $false, $false, $false, $true, 1, $false, "98-21=", 2) | Out-Null
$d.Content.Find.Execute("69-60=", $false, $false, $false, $false, $false, $true, 1, $false, "60-30=", 2) | Out-Null
$d.Content.Find.Execute("19+6=", $false, $false, $false, $false, $false, $true, 1, $false, "54+28=", 2) | Out-Null
$d.Content.Find.Execute("12+50=", $false, $false, $false, $false, $false, $true, 1, $false, "79-15=", 2) | Out-Null
$d.Content.Find.Execute("41-16=", $false, $false, $false, $false, $false, $true, 1, $false, "39-4=", 2) | Out-Null
$d.Content.Find.Execute("44+44=", $false, $false, $false, $false, $false, $true, 1, $false, "24-11=", 2) | Out-Null
$d.Content.Find.Execute("50+7=", $false, $false, $false, $false, $false, $true, 1, $false, "98-83=", 2) | Out-Null
$d.Content.Find.Execute("82+2=", $false, $false, $false, $false, $false, $true, 1, $false, "4+20=", 2) | Out-Null
$d.Content.Find.Execute("89-53=", $false, $false, $false, $false, $false, $true, 1, $false, "71+16=", 2) | Out-Null
$d.Content.Find.Execute("40-25=", $false, $false, $false, $false, $false, $true, 1, $false, "33+55=", 2) | Out-Null
$d.Content.Find.Execute("32+14=", $false, $false, $false, $false, $false, $true, 1, $false, "62+14=", 2) | Out-Null
$d.Content.Find.Execute("66-20=", $false, $false, $false, $false, $false, $true, 1, $false, "74+6=", 2) | Out-Null
$d.Content.Find.Execute("91-41=", $false, $false, $false, $false, $false, $true, 1, $false, "20+1=", 2) | Out-Null
$d.Content.Find.Execute("13-2=", $false, $false, $false, $false, $false, $true, 1, $false, "62+0=", 2) | Out-Null
$d.Content.Find.Execute("41+51=", $false, $false, $false, $false, $false, $true, 1, $false, "55-42=", 2) | Out-Null
$d.Content.Find.Execute("1+6=", $false, $false, $false, $false, $false, $true, 1, $false, "84-32=", 2) | Out-Null
$d.Content.Find.Execute("94-45=", $false, $false, $false, $false, $false, $true, 1, $false, "24-10=", 2) | Out-Null
$d.Content.Find.Execute("92-85=", $false, $false, $false, $false, $false, $true, 1, $false, "52-15=", 2) | Out-Null
$d.Content.Find.Execute("76+23=", $false, $false, $false, $false, $false, $true, 1, $false, "97-18=", 2) | Out-Null
$d.Content.Find.Execute("28+28=", $false, $false, $false, $false, $false, $true, 1, $false, "14+3=", 2) | Out-Null
$d.Content.Find.Execute("81-65=", $false, $false, $false, $false, $false, $true, 1, $false, "92-61=", 2) | Out-Null
$d.Content.Find.Execute("88-79=", $false, $false, $false, $false, $false, $true, 1, $false, "73-41=", 2) | Out-Null
$d.Content.Find.Execute("66-43=", $false, $false, $false, $false, $false, $true, 1, $false, "72-13=", 2) | Out-Null
$d.Content.Find.Execute("13-0=", $false, $false, $false, $false, $false, $true, 1, $false, "66+31=", 2) | Out-Null
$d.Content.Find.Execute("3+12=", $false, $false, $false, $false, $false, $true, 1, $false, "75-58=", 2) | Out-Null
$d.Content.Find.Execute("8+18=", $false, $false, $false, $false, $false, $true, 1, $false, "91-37=", 2) | Out-Null
$d.Content.Find.Execute("89-17=", $false, $false, $false, $false, $false, $true, 1, $false, "28+56=", 2) | Out-Null
$d.Content.Find.Execute("10+38=", $false, $false, $false, $false, $false, $true, 1, $false, "12-6=", 2) | Out-Null
$d.Content.Find.Execute("87-0=", $false, $false, $false, $false, $false, $true, 1, $false, "16-1=", 2) | Out-Null
$d.Content.Find.Execute("45+32=", $false, $false, $false, $false, $false, $true, 1, $false, "18+29=", 2) | Out-Null
$d.Content.Find.Execute("47+0=", $false, $false, $false, $false, $false, $true, 1, $false, "82-51=", 2) | Out-Null
$d.Content.Find.Execute("10+8=", $false, $false, $false, $false, $false, $true, 1, $false, "56-25=", 2) | Out-Null
$d.Content.Find.Execute("83-29=", $false, $false, $false, $false, $false, $true, 1, $false, "3+70=", 2) | Out-Null
$d.Content.Find.Execute("41+1=", $false, $false, $false, $false, $false, $true, 1, $false, "30+27=", 2) | Out-Null
$d.Content.Find.Execute("15+29=", $false, $false, $false, $false, $false, $true, 1, $false, "56+10=", 2) | Out-Null
$d.Content.Find.Execute("25+43=", $false, $false, $false, $false, $false, $true, 1, $false, "30-16=", 2) | Out-Null
$d.Content.Find.Execute("45-42=", $false, $false, $false, $false, $false, $true, 1, $false, "25+4=", 2) | Out-Null
$d.Content.Find.Execute("4+25=", $false, $false, $false, $false, $false, $true, 1, $false, "33-8=", 2) | Out-Null
$d.Content.Find.Execute("80+4=", $false, $false, $false, $false, $false, $true, 1, $false, "17+64=", 2) | Out-Null
$d.Content.Find.Execute("92-31=", $false, $false, $false, $false, $false, $true, 1, $false, "94-59=", 2) | Out-Null
$d.Content.Find.Execute("21+36=", $false, $false, $false, $false, $false, $true, 1, $false, "82-17=", 2) | Out-Null
$d.Content.Find.Execute("12+46=", $false, $false, $false, $false, $false, $true, 1, $false, "26+26=", 2) | Out-Null
$d.Content.Find.Execute("90-29=", $false, $false, $false, $false, $false, $true, 1, $false, "36+43=", 2) | Out-Null
$d.Content.Find.Execute("73-63=", $false, $false, $false, $false, $false, $true, 1, $false, "77-26=", 2) | Out-Null
$d.Content.Find.Execute("19+1=", $false, $false, $false, $false, $false, $true, 1, $false, "59-25=", 2) | Out-Null
$d.Content.Find.Execute("21-0=", $false, $false, $false, $false, $false, $true, 1, $false, "5+75=", 2) | Out-Null
$d.Content.Find.Execute("79-50=", $false, $false, $false, $false, $false, $true, 1, $false, "63+36=", 2) | Out-Null
$d.Content.Find.Execute("68+15=", $false, $false, $false, $false, $false, $true, 1, $false, "18-11=", 2) | Out-Null
$d.Content.Find.Execute("99-26=", $false, $false, $false, $false, $false, $true, 1, $false, "16-3=", 2) | Out-Null
$d.Content.Find.Execute("50+4=", $false, $false, $false, $false, $false, $true, 1, $false, "85-79=", 2) | Out-Null
$d.Content.Find.Execute("46+25=", $false, $false, $false, $false, $false, $true, 1, $false, "66-10=", 2) | Out-Null
$d.Content.Find.Execute("48-9=", $false, $false, $false, $false, $false, $true, 1, $false, "2+78=", 2) | Out-Null
$d.Content.Find.Execute("3+31=", $false, $false, $false, $false, $false, $true, 1, $false, "11+38=", 2) | Out-Null
$d.Content.Find.Execute("85-50=", $false, $false, $false, $false, $false, $true, 1, $false, "17+73=", 2) | Out-Null
$d.Content.Find.Execute("96-37=", $false, $false, $false, $false, $false, $true, 1, $false, "29+20=", 2) | Out-Null
$d.Content.Find.Execute("64-59=", $false, $false, $false, $false, $false, $true, 1, $false, "89-55=", 2) | Out-Null
$d.Content.Find.Execute("38-5=", $false, $false, $false, $false, $false, $true, 1, $false, "85-58=", 2) | Out-Null
$d.Content.Find.Execute("9+36=", $false, $false, $false, $false, $false, $true, 1, $false, "99-59=", 2) | Out-Null
$d.Content.Find.Execute("4+14=", $false, $false, $false, $false, $false, $true, 1, $false, "24+30=", 2) | Out-Null
$d.Content.Find.Execute("7+16=", $false, $false, $false, $false, $false, $true, 1, $false, "56-42=", 2) | Out-Null
$d.Content.Find.Execute("97-65=", $false, $false, $false, $false, $false, $true, 1, $false, "23-1=", 2) | Out-Null
$d.Content.Find.Execute("16+13=", $false, $false, $false, $false, $false, $true, 1, $false, "12+68=", 2) | Out-Null
$d.Content.Find.Execute("52-36=", $false, $false, $false, $false, $false, $true, 1, $false, "59-1=", 2) | Out-Null
$d.Content.Find.Execute("12-9=", $false, $false, $false, $false, $false, $true, 1, $false, "95-18=", 2) | Out-Null
$d.Content.Find.Execute("24+2=", $false, $false, $false, $false, $false, $true, 1, $false, "59+15=", 2) | Out-Null
$d.Content.Find.Execute("90-50=", $false, $false, $false, $false, $false, $true, 1, $false, "63-40=", 2) | Out-Null
$d.Content.Find.Execute("45-17=", $false, $false, $false, $false, $false, $true, 1, $false, "62-11=", 2) | Out-Null
$d.Content.Find.Execute("36+37=", $false, $false, $false, $false, $false, $true, 1, $false, "64-15=", 2) | Out-Null
$d.Content.Find.Execute("20+3=", $false, $false, $false, $false, $false, $true, 1, $false, "71+11=", 2) | Out-Null
$d.Content.Find.Execute("55-50=", $false, $false, $false, $false, $false, $true, 1, $false, "30-9=", 2) | Out-Null
$d.Content.Find.Execute("38+61=", $false, $false, $false, $false, $false, $true, 1, $false, "76-55=", 2) | Out-Null
$d.Content.Find.Execute("90-87=", $false, $false, $false, $false, $false, $true, 1, $false, "82-42=", 2) | Out-Null
$d.Content.Find.Execute("67-55=", $false, $false, $false, $false, $false, $true, 1, $false, "22+62=", 2) | Out-Null
